# Insert a new weekly price record as row 58, pushing all existing rows
# (58..175) down by one (to 59..176). This mirrors the commit
# "Fruta / hortaliza, semanal": a new weekly observation is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 58..175 down to 59..176, leaving a blank row 58.
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new record.
$ws.Range('A58').Value = 4
$ws.Range('B58').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C58').Value = 'Los Lagos'
$ws.Range('D58').Value = 44469
$ws.Range('E58').Value = 10
$ws.Range('F58').Value = 100112045
$ws.Range('G58').Value = 'Zapallo'
$ws.Range('H58').Value = 'Paine'
$ws.Range('I58').Value = '1a (guarda)'
$ws.Range('J58').Value = 500
$ws.Range('K58').Value = 600
$ws.Range('L58').Value = 600
$ws.Range('M58').Value = 600
$ws.Range('N58').Value = '$/kilo (volumen en unidades)'
$ws.Range('O58').Value = 'Región Metropolitana'
$ws.Range('P58').Value = 600
$ws.Range('Q58').Value = 1
$ws.Range('R58').Value = 'Hortaliza'
